# Add data for 2022-09-17 (one additional day of carjacking reports
# rolled into the "through" month-to-date column, plus assorted
# historical September backfills for prior years).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "through" date in the sheet tab name and the column header ---
$ws.Name = "Through 2022-09-09"
$ws.Range("B1").Value = "September 2022 (through September 09)"

# --- Update existing counts ---
$ws.Range("K3").Value = 6     # Garfield Park, September 2021
$ws.Range("BD18").Value = 2   # Ashburn, September 2016
$ws.Range("K21").Value = 2    # Logan Square, September 2021
$ws.Range("B24").Value = 3    # South Shore, current month-to-date
$ws.Range("AL24").Value = 1   # South Shore, September 2018

# --- New counts (previously-empty cells) ---
$ws.Range("AU14").Value = 1   # Roseland, September 2017
$ws.Range("AU27").Value = 1   # Belmont Cragin, September 2017
$ws.Range("K31").Value = 1    # Streeterville, September 2021
$ws.Range("AL31").Value = 1   # Streeterville, September 2018
$ws.Range("K32").Value = 1    # Ukrainian Village, September 2021
$ws.Range("T33").Value = 1    # United Center, September 2020
$ws.Range("AL40").Value = 1   # Pullman, September 2018
$ws.Range("AU43").Value = 1   # Avondale, September 2017
$ws.Range("K50").Value = 1    # Grand Boulevard, September 2021
$ws.Range("BD50").Value = 1   # Grand Boulevard, September 2016
$ws.Range("K69").Value = 1    # Gold Coast, September 2021
$ws.Range("AC78").Value = 1   # Lincoln Park, September 2019
